# repull data, push all data, mean calculation
# Updates column F ("dSF") values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -9
    4  = 1
    6  = 0
    8  = -1
    12 = 4
    16 = 1
    17 = -1
    20 = -5
    23 = -4
    24 = -4
    25 = -5
    28 = 1
    36 = -5
    37 = 5
    38 = -1
    39 = 10
    41 = -4
    42 = -13
    43 = -7
    44 = -3
    45 = -4
    49 = -6
    50 = 5
    51 = -2
    53 = -3
    54 = 4
    55 = -4
    60 = -1
    65 = 1
    68 = -12
    69 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
